# Update the LR-pairs sheet with the new TPM-derived values.
#
# Summary of the change (per the commit "update scripts wuth new tpm"):
#  - The "Sending cluster" / "Target cluster" category used for the 3rd
#    sending cluster and the single target cluster changes from "MuSCs"
#    to a richer set of clusters ("Resolving-Mac", "Inflammatory-Mac",
#    "MuSCs", ...), so every data row's cluster labels and the 20 numeric
#    metric columns (E..T) are refreshed with newly computed values.
#  - The data block shrinks from 16 rows (rows 2-17) to 15 rows (rows 2-16).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-obsolete last data row (old row 17) first, while row
# indices still match the original layout.
$ws.Rows.Item(17).Delete()

# Column order for each record: A..T
#   A = Sending cluster, B = Ligand symbol, C = Receptor symbol,
#   D = Target cluster, E..T = numeric metrics
$data = @()
$data += ,@("ECs", "Colq", "Musk", "ECs", 1, 0.3333333333333333, 0.0313, 0.0939, 0.2134634270541592, 0.2134634270541592, 1, 0.3333333333333333, 0.1113626666666667, 0.334088, 0.005543586609007977, 0.00750796989180202, 0.003485651466666667, 0.0313708632, 0.001183352995730388, 0.001602676983323504)
$data += ,@("ECs", "Colq", "Musk", "FAPs", 1, 0.3333333333333333, 0.0313, 0.0939, 0.2134634270541592, 0.2134634270541592, 3, 1, 4.120856, 12.362568, 0.2051344748022992, 0.2778243706129975, 0.1289827928, 1.1608451352, 0.04378870799825385, 0.0593053422702153)
$data += ,@("ECs", "Colq", "Musk", "Inflammatory-Mac", 1, 0.3333333333333333, 0.0313, 0.0939, 0.2134634270541592, 0.2134634270541592, 1, 0.3333333333333333, 0.01920333333333333, 0.05761, 0.0009559338394223962, 0.00129467130057564, 0.0006010643333333334, 0.005409579, 0.0002040569134001451, 0.0002763649727295416)
$data += ,@("ECs", "Colq", "Musk", "MuSCs", 1, 0.3333333333333333, 0.0313, 0.0939, 0.2134634270541592, 0.2134634270541592, 2, 1, 15.767896, 31.535792, 0.7849192169532917, 0.7087048228315025, 0.4935351448, 2.9612108688, 0.1675515460115168, 0.1512825602514233)
$data += ,@("ECs", "Colq", "Musk", "Resolving-Mac", 1, 0.3333333333333333, 0.0313, 0.0939, 0.2134634270541592, 0.2134634270541592, 2, 0.6666666666666666, 0.06924100000000001, 0.207723, 0.003446787795978796, 0.004668165363122264, 0.0021672433, 0.0195051897, 0.000735763135258086, 0.0009964825764676023)
$data += ,@("FAPs", "Colq", "Musk", "ECs", 3, 1, 0.05318533333333333, 0.159556, 0.3627196013530717, 0.3627196013530717, 1, 0.3333333333333333, 0.1113626666666667, 0.334088, 0.005543586609007977, 0.00750796989180202, 0.005922860547555556, 0.053305744928, 0.0020107675248856, 0.002723287846125294)
$data += ,@("FAPs", "Colq", "Musk", "FAPs", 3, 1, 0.05318533333333333, 0.159556, 0.3627196013530717, 0.3627196013530717, 3, 1, 4.120856, 12.362568, 0.2051344748022992, 0.2778243706129975, 0.2191690999786667, 1.972521899808, 0.07440629492406167, 0.1007723449549145)
$data += ,@("FAPs", "Colq", "Musk", "Inflammatory-Mac", 3, 1, 0.05318533333333333, 0.159556, 0.3627196013530717, 0.3627196013530717, 1, 0.3333333333333333, 0.01920333333333333, 0.05761, 0.0009559338394223962, 0.00129467130057564, 0.001021335684444444, 0.009192021160000001, 0.0003467359411552028, 0.000469602658028059)
$data += ,@("FAPs", "Colq", "Musk", "MuSCs", 3, 1, 0.05318533333333333, 0.159556, 0.3627196013530717, 0.3627196013530717, 2, 1, 15.767896, 31.535792, 0.7849192169532917, 0.7087048228315025, 0.8386208047253334, 5.031724828352, 0.2847055854676632, 0.2570611308144419)
$data += ,@("FAPs", "Colq", "Musk", "Resolving-Mac", 3, 1, 0.05318533333333333, 0.159556, 0.3627196013530717, 0.3627196013530717, 2, 0.6666666666666666, 0.06924100000000001, 0.207723, 0.003446787795978796, 0.004668165363122264, 0.003682605665333334, 0.03314345098800001, 0.001250217495306061, 0.001693235079561925)
$data += ,@("Resolving-Mac", "Colq", "Musk", "ECs", 1, 0.3333333333333333, 0.062144, 0.186432, 0.423816971592769, 0.423816971592769, 1, 0.3333333333333333, 0.1113626666666667, 0.334088, 0.005543586609007977, 0.00750796989180202, 0.006920521557333333, 0.06228469401599999, 0.002349466088391988, 0.003182005062353222)
$data += ,@("Resolving-Mac", "Colq", "Musk", "FAPs", 1, 0.3333333333333333, 0.062144, 0.186432, 0.423816971592769, 0.423816971592769, 3, 1, 4.120856, 12.362568, 0.2051344748022992, 0.2778243706129975, 0.256086475264, 2.304778277376, 0.08693947187998362, 0.1177466833878677)
$data += ,@("Resolving-Mac", "Colq", "Musk", "Inflammatory-Mac", 1, 0.3333333333333333, 0.062144, 0.186432, 0.423816971592769, 0.423816971592769, 1, 0.3333333333333333, 0.01920333333333333, 0.05761, 0.0009559338394223962, 0.00129467130057564, 0.001193371946666667, 0.01074034752, 0.0004051409848670483, 0.0005487036698180394)
$data += ,@("Resolving-Mac", "Colq", "Musk", "MuSCs", 1, 0.3333333333333333, 0.062144, 0.186432, 0.423816971592769, 0.423816971592769, 2, 1, 15.767896, 31.535792, 0.7849192169532917, 0.7087048228315025, 0.979880129024, 5.879280774144, 0.3326620854741117, 0.3003611317656373)
$data += ,@("Resolving-Mac", "Colq", "Musk", "Resolving-Mac", 1, 0.3333333333333333, 0.062144, 0.186432, 0.423816971592769, 0.423816971592769, 2, 0.6666666666666666, 0.06924100000000001, 0.207723, 0.003446787795978796, 0.004668165363122264, 0.004302912704000001, 0.038726214336, 0.001460807165414649, 0.001978447707092737)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $rowValues = $data[$i]
    $r = $startRow + $i
    for ($c = 1; $c -le $rowValues.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowValues[$c - 1]
    }
}
